# Tutorial-04 "Network IP Assignments.xlsx" — clarify the workflow for
# optimal IP assignment:
#   1. Rename the summary sheet "Sheet1" -> "IP Ranges" so its purpose is
#      obvious in the tab strip.
#   2. Make "IP Ranges" the active/selected sheet on open (it was
#      "Network N2" before), since that's now the natural starting point
#      of the workflow.

$wb = $excel.ActiveWorkbook

$ipRanges = $wb.Worksheets.Item("Sheet1")
$ipRanges.Name = "IP Ranges"

$ipRanges.Activate()
